# Journal_Travail_Francois.xlsx - add week 12 entries (rows 67-70) + pieChart /
# progress-bar related journal lines, matching the author's commit
# "pieChart + maj jdt francois".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Copy the visual formatting (borders/fonts/number formats) from an
#    existing 4-row "week" block (rows 16:19 -> week 4) onto the currently
#    blank placeholder rows 67:70, which will become week 12. This preserves
#    the thin/medium border pattern + merged "Semaines"/"Total heures" look
#    used by every other week group instead of the generic blank-row style.
# ---------------------------------------------------------------------------
$ws.Range("A16:E19").Copy()
$ws.Range("A67").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Fill in the new journal entries for week 12.
# ---------------------------------------------------------------------------
$ws.Range("A67").Value = 43204
$ws.Range("B67").Value = "Ajout des images à tout les bouton (add/supp/retour/edit/setting)"
$ws.Range("C67").Value = 0.75

$ws.Range("A68").Value = 43204
$ws.Range("B68").Value = "Calcule des dépense sur plusieurs catégorie"
$ws.Range("C68").Value = 0.75

$ws.Range("A69").Value = 43204
$ws.Range("B69").Value = "Progresse bar sur les budget"
$ws.Range("C69").Value = 0.75

$ws.Range("A70").Value = 43234
$ws.Range("B70").Value = "Graphique en camembert budget"
$ws.Range("C70").Value = 0.75

# Week number + weekly total (merged across the 4-row block, like the other
# week groups).
$ws.Range("D67").Value = 12
$ws.Range("E67").Formula = "=SUM(C67:C70)"

# ---------------------------------------------------------------------------
# 3) Merge the week-number / weekly-total columns over the new block.
# ---------------------------------------------------------------------------
$ws.Range("D67:D70").Merge()
$ws.Range("E67:E70").Merge()

# ---------------------------------------------------------------------------
# 4) Update the sheet view to match where the author left the cursor.
# ---------------------------------------------------------------------------
$win = $excel.Windows.Item(1)
$win.ScrollRow = 61
$win.ScrollColumn = 1
$ws.Range("F70").Select()
